# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from 45224 (2023-10-25) to 45233 (2023-11-03), keeping the existing
# date number formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value2 = 45233
    }
}
